$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: remove the "Meta description" paragraph that currently sits
# right after the title heading near the top of the document.
# ---------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("Meta description")
if ($found) {
    $metaPara = $rng.Paragraphs(1)
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------
# Step 2: at the end of the document, insert a new bold heading
# paragraph ("Play Bondi Break Slot for Free - Review & Demo Game")
# right before the final (italic) paragraph, and replace that final
# paragraph's text with the former meta-description copy (kept italic).
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)

$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Bondi Break Slot for Free - Review &amp; Demo Game</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Experience the beach with Bondi Break slot. Play for free with engaging gameplay, potential wins of up to 6,250x your stake and a lively beach-themed design.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$lastPara.Range.InsertXML($xml)

# InsertXML leaves the document's mandatory trailing paragraph mark
# behind as an extra empty paragraph - merge it back into the
# paragraph above it so the body ends cleanly on our new text.
$trailing = $d.Paragraphs($d.Paragraphs.Count)
if ($trailing.Range.Text.Trim() -eq "") {
    $prev = $d.Paragraphs($d.Paragraphs.Count - 1)
    $mergeRange = $d.Range($prev.Range.End - 1, $prev.Range.End)
    $mergeRange.Delete()
}
